# Realign CryCompanywiseStockReport rows whose Material Code / Rate / Qty / Value
# columns were shifted by one position within their product group (each group
# repeats a description across consecutive rows), and bump the report
# From/To date header (I1, K1) from 01-10-2025 to 02-10-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "02-10-2025 00:00:00"
$ws.Range("K1").Value = "02-10-2025 00:00:00"

# Row 142
$ws.Cells.Item(142, 2).Value = 48654
$ws.Cells.Item(142, 3).Value = "CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms"
$ws.Cells.Item(142, 4).Value = 32.02
$ws.Cells.Item(142, 5).Value = 38.26
$ws.Cells.Item(142, 6).Value = -1
$ws.Cells.Item(142, 7).Value = -32.02

# Row 143
$ws.Cells.Item(143, 2).Value = 63902
$ws.Cells.Item(143, 3).Value = "CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms"
$ws.Cells.Item(143, 4).Value = 32.02
$ws.Cells.Item(143, 5).Value = 34.04
$ws.Cells.Item(143, 6).Value = 2
$ws.Cells.Item(143, 7).Value = 64.04000000000001

# Row 154
$ws.Cells.Item(154, 2).Value = 64350
$ws.Cells.Item(154, 3).Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Cells.Item(154, 4).Value = 66.44
$ws.Cells.Item(154, 5).Value = 70.63
$ws.Cells.Item(154, 6).Value = 101
$ws.Cells.Item(154, 7).Value = 6710.44

# Row 155
$ws.Cells.Item(155, 2).Value = 57756
$ws.Cells.Item(155, 3).Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Cells.Item(155, 4).Value = 66.44
$ws.Cells.Item(155, 5).Value = 79.37
$ws.Cells.Item(155, 6).Value = -100
$ws.Cells.Item(155, 7).Value = -6644

# Row 156
$ws.Cells.Item(156, 2).Value = 53925
$ws.Cells.Item(156, 3).Value = "COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush"
$ws.Cells.Item(156, 4).Value = 66.44
$ws.Cells.Item(156, 5).Value = 79.37
$ws.Cells.Item(156, 6).Value = 1
$ws.Cells.Item(156, 7).Value = 66.44

# Row 176
$ws.Cells.Item(176, 2).Value = 57552
$ws.Cells.Item(176, 3).Value = "DAB-Real Activ Coconut Water Tetra 1000ml"
$ws.Cells.Item(176, 4).Value = 120.69
$ws.Cells.Item(176, 5).Value = 136.86
$ws.Cells.Item(176, 6).Value = -5
$ws.Cells.Item(176, 7).Value = -603.45

# Row 177
$ws.Cells.Item(177, 2).Value = 64329
$ws.Cells.Item(177, 3).Value = "DAB-Real Activ Coconut Water Tetra 1000ml"
$ws.Cells.Item(177, 4).Value = 120.69
$ws.Cells.Item(177, 5).Value = 128.32
$ws.Cells.Item(177, 6).Value = 6
$ws.Cells.Item(177, 7).Value = 724.14

# Row 256
$ws.Cells.Item(256, 2).Value = 48719
$ws.Cells.Item(256, 3).Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Cells.Item(256, 4).Value = 295.75
$ws.Cells.Item(256, 5).Value = 353.35
$ws.Cells.Item(256, 6).Value = -81
$ws.Cells.Item(256, 7).Value = -23955.75

# Row 257
$ws.Cells.Item(257, 2).Value = 64979
$ws.Cells.Item(257, 3).Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Cells.Item(257, 4).Value = 295.75
$ws.Cells.Item(257, 5).Value = 314.41
$ws.Cells.Item(257, 6).Value = 82
$ws.Cells.Item(257, 7).Value = 24251.5

# Row 305
$ws.Cells.Item(305, 2).Value = 57854
$ws.Cells.Item(305, 3).Value = "HUL-3Roses Dust [C] 500G Relaunch"
$ws.Cells.Item(305, 4).Value = 305.84
$ws.Cells.Item(305, 5).Value = 325.16
$ws.Cells.Item(305, 6).Value = 2
$ws.Cells.Item(305, 7).Value = 611.6799999999999

# Row 306
$ws.Cells.Item(306, 2).Value = 62997
$ws.Cells.Item(306, 3).Value = "HUL-3Roses Dust [C] 500G Relaunch"
$ws.Cells.Item(306, 4).Value = 305.84
$ws.Cells.Item(306, 5).Value = 325.16
$ws.Cells.Item(306, 6).Value = 72
$ws.Cells.Item(306, 7).Value = 22020.48

# Row 308
$ws.Cells.Item(308, 2).Value = 57077
$ws.Cells.Item(308, 3).Value = "HUL-Bru Inst Poly 50g"
$ws.Cells.Item(308, 4).Value = 93.08
$ws.Cells.Item(308, 5).Value = 111.2
$ws.Cells.Item(308, 6).Value = 1
$ws.Cells.Item(308, 7).Value = 93.08

# Row 309
$ws.Cells.Item(309, 2).Value = 61610
$ws.Cells.Item(309, 3).Value = "HUL-Bru Inst Poly 50g"
$ws.Cells.Item(309, 4).Value = 102.71
$ws.Cells.Item(309, 5).Value = 122.71
$ws.Cells.Item(309, 6).Value = -58
$ws.Cells.Item(309, 7).Value = -5957.18

# Row 310
$ws.Cells.Item(310, 2).Value = 63565
$ws.Cells.Item(310, 3).Value = "HUL-Bru Inst Poly 50g"
$ws.Cells.Item(310, 4).Value = 102.71
$ws.Cells.Item(310, 5).Value = 109.19
$ws.Cells.Item(310, 6).Value = 60
$ws.Cells.Item(310, 7).Value = 6162.6

# Row 338
$ws.Cells.Item(338, 2).Value = 63520
$ws.Cells.Item(338, 3).Value = "HUL-Kissan nango jam 490g"
$ws.Cells.Item(338, 4).Value = 144.28
$ws.Cells.Item(338, 5).Value = 153.4
$ws.Cells.Item(338, 6).Value = 97
$ws.Cells.Item(338, 7).Value = 13995.16

# Row 339
$ws.Cells.Item(339, 2).Value = 55373
$ws.Cells.Item(339, 3).Value = "HUL-Kissan nango jam 490g"
$ws.Cells.Item(339, 4).Value = 144.28
$ws.Cells.Item(339, 5).Value = 163.62
$ws.Cells.Item(339, 6).Value = -94
$ws.Cells.Item(339, 7).Value = -13562.32

# Row 343
$ws.Cells.Item(343, 2).Value = 57802
$ws.Cells.Item(343, 3).Value = "HUL-Kissan Pineapple Jam 500G"
$ws.Cells.Item(343, 4).Value = 143.48
$ws.Cells.Item(343, 5).Value = 162.71
$ws.Cells.Item(343, 6).Value = -79
$ws.Cells.Item(343, 7).Value = -11334.92

# Row 344
$ws.Cells.Item(344, 2).Value = 63571
$ws.Cells.Item(344, 3).Value = "HUL-Kissan Pineapple Jam 500G"
$ws.Cells.Item(344, 4).Value = 143.48
$ws.Cells.Item(344, 5).Value = 152.53
$ws.Cells.Item(344, 6).Value = 29
$ws.Cells.Item(344, 7).Value = 4160.92

# Row 347
$ws.Cells.Item(347, 2).Value = 63510
$ws.Cells.Item(347, 3).Value = "HUL-knorr schezwan 200g pch"
$ws.Cells.Item(347, 4).Value = 47.64
$ws.Cells.Item(347, 5).Value = 50.66
$ws.Cells.Item(347, 6).Value = 167
$ws.Cells.Item(347, 7).Value = 7955.88

# Row 348
$ws.Cells.Item(348, 2).Value = 55356
$ws.Cells.Item(348, 3).Value = "HUL-knorr schezwan 200g pch"
$ws.Cells.Item(348, 4).Value = 47.64
$ws.Cells.Item(348, 5).Value = 54.04
$ws.Cells.Item(348, 6).Value = -158
$ws.Cells.Item(348, 7).Value = -7527.12

# Row 381
$ws.Cells.Item(381, 2).Value = 57817
$ws.Cells.Item(381, 3).Value = "HUL-Rap Refresh Bolt 1Kg"
$ws.Cells.Item(381, 4).Value = 79.81
$ws.Cells.Item(381, 5).Value = 95.34999999999999
$ws.Cells.Item(381, 6).Value = 3
$ws.Cells.Item(381, 7).Value = 239.43

# Row 382
$ws.Cells.Item(382, 2).Value = 62865
$ws.Cells.Item(382, 3).Value = "HUL-Rap Refresh Bolt 1Kg"
$ws.Cells.Item(382, 4).Value = 79.81
$ws.Cells.Item(382, 5).Value = 95.34999999999999
$ws.Cells.Item(382, 6).Value = 151
$ws.Cells.Item(382, 7).Value = 12051.31

# Row 411
$ws.Cells.Item(411, 2).Value = 57856
$ws.Cells.Item(411, 3).Value = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Cells.Item(411, 4).Value = 171.33
$ws.Cells.Item(411, 5).Value = 204.69
$ws.Cells.Item(411, 6).Value = 2
$ws.Cells.Item(411, 7).Value = 342.66

# Row 412
$ws.Cells.Item(412, 2).Value = 63007
$ws.Cells.Item(412, 3).Value = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Cells.Item(412, 4).Value = 171.33
$ws.Cells.Item(412, 5).Value = 204.69
$ws.Cells.Item(412, 6).Value = 984
$ws.Cells.Item(412, 7).Value = 168588.72

# Row 413
$ws.Cells.Item(413, 2).Value = 57857
$ws.Cells.Item(413, 3).Value = "HUL-Surf Exl Mtc Liq Tl 1 Ltr Cp"
$ws.Cells.Item(413, 4).Value = 151.17
$ws.Cells.Item(413, 5).Value = 180.62
$ws.Cells.Item(413, 6).Value = 3
$ws.Cells.Item(413, 7).Value = 453.51

# Row 414
$ws.Cells.Item(414, 2).Value = 63008
$ws.Cells.Item(414, 3).Value = "HUL-Surf Exl Mtc Liq Tl 1 Ltr Cp"
$ws.Cells.Item(414, 4).Value = 151.17
$ws.Cells.Item(414, 5).Value = 180.62
$ws.Cells.Item(414, 6).Value = 504
$ws.Cells.Item(414, 7).Value = 76189.67999999999

# Row 423
$ws.Cells.Item(423, 2).Value = 63102
$ws.Cells.Item(423, 3).Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Cells.Item(423, 4).Value = 59.47
$ws.Cells.Item(423, 5).Value = 71.05
$ws.Cells.Item(423, 6).Value = 36
$ws.Cells.Item(423, 7).Value = 2140.92

# Row 424
$ws.Cells.Item(424, 2).Value = 53082
$ws.Cells.Item(424, 3).Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Cells.Item(424, 4).Value = 59.47
$ws.Cells.Item(424, 5).Value = 71.05
$ws.Cells.Item(424, 6).Value = 1
$ws.Cells.Item(424, 7).Value = 59.47

# Row 449
$ws.Cells.Item(449, 2).Value = 31930
$ws.Cells.Item(449, 3).Value = "JLM-MBD Shiny Toothbrush Safari"
$ws.Cells.Item(449, 4).Value = 22.42
$ws.Cells.Item(449, 5).Value = 26.8
$ws.Cells.Item(449, 6).Value = -62
$ws.Cells.Item(449, 7).Value = -1390.04

# Row 450
$ws.Cells.Item(450, 2).Value = 63681
$ws.Cells.Item(450, 3).Value = "JLM-MBD Shiny Toothbrush Safari"
$ws.Cells.Item(450, 4).Value = 22.42
$ws.Cells.Item(450, 5).Value = 23.84
$ws.Cells.Item(450, 6).Value = 65
$ws.Cells.Item(450, 7).Value = 1457.3

# Row 528
$ws.Cells.Item(528, 2).Value = 58047
$ws.Cells.Item(528, 3).Value = "KUS-Floor Wiper"
$ws.Cells.Item(528, 4).Value = 105.54
$ws.Cells.Item(528, 5).Value = 126.1
$ws.Cells.Item(528, 6).Value = 54
$ws.Cells.Item(528, 7).Value = 5699.16

# Row 529
$ws.Cells.Item(529, 2).Value = 47097
$ws.Cells.Item(529, 3).Value = "KUS-Floor Wiper"
$ws.Cells.Item(529, 4).Value = 112.28
$ws.Cells.Item(529, 5).Value = 134.16
$ws.Cells.Item(529, 6).Value = 15
$ws.Cells.Item(529, 7).Value = 1684.2

# Row 575
$ws.Cells.Item(575, 2).Value = 65066
$ws.Cells.Item(575, 3).Value = "CRE-Butter cremfills 100gm"
$ws.Cells.Item(575, 4).Value = 12.81
$ws.Cells.Item(575, 5).Value = 13.61
$ws.Cells.Item(575, 6).Value = 313
$ws.Cells.Item(575, 7).Value = 4009.53

# Row 576
$ws.Cells.Item(576, 2).Value = 53263
$ws.Cells.Item(576, 3).Value = "CRE-Butter cremfills 100gm"
$ws.Cells.Item(576, 4).Value = 12.81
$ws.Cells.Item(576, 5).Value = 15.29
$ws.Cells.Item(576, 6).Value = -309
$ws.Cells.Item(576, 7).Value = -3958.29

# Row 582
$ws.Cells.Item(582, 2).Value = 64922
$ws.Cells.Item(582, 3).Value = "CRE-Cremica Golden Bytes Rich Butter 200Gm"
$ws.Cells.Item(582, 4).Value = 19.73
$ws.Cells.Item(582, 5).Value = 20.98
$ws.Cells.Item(582, 6).Value = 207
$ws.Cells.Item(582, 7).Value = 4084.11

# Row 583
$ws.Cells.Item(583, 2).Value = 45706
$ws.Cells.Item(583, 3).Value = "CRE-Cremica Golden Bytes Rich Butter 200Gm"
$ws.Cells.Item(583, 4).Value = 19.73
$ws.Cells.Item(583, 5).Value = 23.58
$ws.Cells.Item(583, 6).Value = -202
$ws.Cells.Item(583, 7).Value = -3985.46

# Row 585
$ws.Cells.Item(585, 2).Value = 45718
$ws.Cells.Item(585, 3).Value = "CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm"
$ws.Cells.Item(585, 4).Value = 16.22
$ws.Cells.Item(585, 5).Value = 19.38
$ws.Cells.Item(585, 6).Value = -294
$ws.Cells.Item(585, 7).Value = -4768.68

# Row 586
$ws.Cells.Item(586, 2).Value = 64927
$ws.Cells.Item(586, 3).Value = "CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm"
$ws.Cells.Item(586, 4).Value = 16.22
$ws.Cells.Item(586, 5).Value = 17.26
$ws.Cells.Item(586, 6).Value = 295
$ws.Cells.Item(586, 7).Value = 4784.9

# Row 596
$ws.Cells.Item(596, 2).Value = 53595
$ws.Cells.Item(596, 3).Value = "CRE-Kaju khz cookies 100 gm"
$ws.Cells.Item(596, 4).Value = 14.73
$ws.Cells.Item(596, 5).Value = 17.61
$ws.Cells.Item(596, 6).Value = -335
$ws.Cells.Item(596, 7).Value = -4934.55

# Row 597
$ws.Cells.Item(597, 2).Value = 65067
$ws.Cells.Item(597, 3).Value = "CRE-Kaju khz cookies 100 gm"
$ws.Cells.Item(597, 4).Value = 14.73
$ws.Cells.Item(597, 5).Value = 15.65
$ws.Cells.Item(597, 6).Value = 338
$ws.Cells.Item(597, 7).Value = 4978.74

# Row 679
$ws.Cells.Item(679, 2).Value = 53319
$ws.Cells.Item(679, 3).Value = "PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)"
$ws.Cells.Item(679, 4).Value = 273.92
$ws.Cells.Item(679, 5).Value = 310.64
$ws.Cells.Item(679, 6).Value = -6
$ws.Cells.Item(679, 7).Value = -1643.52

# Row 680
$ws.Cells.Item(680, 2).Value = 64810
$ws.Cells.Item(680, 3).Value = "PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)"
$ws.Cells.Item(680, 4).Value = 273.92
$ws.Cells.Item(680, 5).Value = 291.22
$ws.Cells.Item(680, 6).Value = 7
$ws.Cells.Item(680, 7).Value = 1917.44

# Row 712
$ws.Cells.Item(712, 2).Value = 60022
$ws.Cells.Item(712, 3).Value = "Rasna Nagpur Orange (32 Glass)"
$ws.Cells.Item(712, 4).Value = 32.83
$ws.Cells.Item(712, 5).Value = 37.22
$ws.Cells.Item(712, 6).Value = -113
$ws.Cells.Item(712, 7).Value = -3709.79

# Row 713
$ws.Cells.Item(713, 2).Value = 64830
$ws.Cells.Item(713, 3).Value = "Rasna Nagpur Orange (32 Glass)"
$ws.Cells.Item(713, 4).Value = 32.83
$ws.Cells.Item(713, 5).Value = 34.9
$ws.Cells.Item(713, 6).Value = 117
$ws.Cells.Item(713, 7).Value = 3841.11

# Row 864
$ws.Cells.Item(864, 2).Value = 54751
$ws.Cells.Item(864, 3).Value = "Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm"
$ws.Cells.Item(864, 4).Value = 40.87
$ws.Cells.Item(864, 5).Value = 46.34
$ws.Cells.Item(864, 6).Value = -19
$ws.Cells.Item(864, 7).Value = -776.53

# Row 865
$ws.Cells.Item(865, 2).Value = 65079
$ws.Cells.Item(865, 3).Value = "Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm"
$ws.Cells.Item(865, 4).Value = 40.87
$ws.Cells.Item(865, 5).Value = 43.44
$ws.Cells.Item(865, 6).Value = 21
$ws.Cells.Item(865, 7).Value = 858.27
